$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts existing A:E -> B:F
$ws.Columns.Item(1).Insert()

# Give the new header cell (A1) the same style as its neighbour (bold/bordered
# header style), then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1").Value = "ID"

# Populate the new ID column for each data row.
$ids = @(
    @{Row=2;  Id="Hb 2"},
    @{Row=3;  Id="Hb 3"},
    @{Row=4;  Id="S 24"},
    @{Row=5;  Id="S 28"},
    @{Row=6;  Id="Hb 107"},
    @{Row=7;  Id="Hb 66"},
    @{Row=8;  Id="Hb 69"},
    @{Row=9;  Id="Hb 95"},
    @{Row=10; Id="Hb 99"},
    @{Row=11; Id="Hb 92"},
    @{Row=12; Id="Hb 40"},
    @{Row=13; Id="Hb 41"},
    @{Row=14; Id="S 11"},
    @{Row=15; Id="Hb 57"},
    @{Row=16; Id="S 21"},
    @{Row=17; Id="S 22"},
    @{Row=18; Id="S 3"},
    @{Row=19; Id="S 4"},
    @{Row=20; Id="S 5"},
    @{Row=21; Id="Hb 74"},
    @{Row=22; Id="Hb 79"},
    @{Row=23; Id="Hb 32"},
    @{Row=24; Id="S 15"},
    @{Row=25; Id="S 16"}
)

foreach ($entry in $ids) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Id
}

Write-Host "Inserted ID column and populated labels"
